$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Existing-cell edits -----------------------------------------------
# H14: "poster" -> "oral"
$ws.Range("H14").Value = "oral"

# H26: "poster" -> "oral"
$ws.Range("H26").Value = "oral"

# --- New rows 46-47: paper title, then author info ----------------------
$ws.Range("G46").Value = "Commuting Conjugacy Class Graph of The Generalized"
$ws.Range("G47").Value = "Dicyclic Group"

$ws.Range("A46").Value = "محمد علی"
$ws.Range("B46").Value = "سلحشور"
$ws.Range("C46").Value = "آزاد سواد کوه"
$ws.Range("D46").Value = "گروه"
$ws.Range("E46").Value = "MA.Salahshour@iau.ac.ir"
$ws.Range("F46").Value = 550000
$ws.Range("H46").Value = "oral"

# C28: "مجتمع آموزش عالی بم" -> "پیام نور " (trailing space)
$ws.Range("C28").Value = "پیام نور "

# --- New row 48 -----------------------------------------------------------
$ws.Range("B48").Value = "قدیمی"
$ws.Range("A48").Value = "کریم"
$ws.Range("F48").Value = 550000
$ws.Range("G48").Value = "A way to construction strong d-algebras"
$ws.Range("H48").Value = "oral"

# --- View state: selection moves to F48 ----------------------------------
[void]$ws.Range("F48").Select()
